$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: "ID" -> "Id"
$ws.Range("A1").Value = "Id"

# Update header text: "Duration (h)" -> "Duration"
$ws.Range("D1").Value = "Duration"

# Move the active selection to A2 (matches the new cursor position in the saved file)
$ws.Range("A2").Select()
